$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on all target cells so values like "1.001" or
# "30.035.12" are stored as literal text (matching the inline string cells
# in the original workbook) rather than being auto-converted to numbers.
$targets = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "B12", "C12", "D12", "E12", "B13", "C13", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "E31", "E32", "D33", "D34", "E34", "D35", "E35", "D36", "E36", "E37", "D38", "E38", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "E51")
foreach ($ref in $targets) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '30.035.12'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '1.859.07'
$ws.Range('E3').Value = '  -2.85%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '232.26'
$ws.Range('E5').Value = '  -2.92%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  -2.56%  '
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('D9').Value = '0.06526'
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('D10').Value = '19.97'
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('D11').Value = '0.07741'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.900.92'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = '95.94'
$ws.Range('E13').Value = '  -6.70%  '
$ws.Range('D14').Value = '5.040'
$ws.Range('E14').Value = '  -3.40%  '
$ws.Range('D15').Value = '0.6650'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '280.01'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '30.063.50'
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').Value = '2.120.77'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').Value = '5.310'
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '0.000007193'
$ws.Range('E23').Value = '  -3.80%  '
$ws.Range('D24').Value = '6.111'
$ws.Range('E24').Value = '  -3.25%  '
$ws.Range('D25').Value = '166.56'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').Value = '9.240'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').Value = '18.86'
$ws.Range('E27').Value = '  -2.75%  '
$ws.Range('D28').Value = '1.926'
$ws.Range('E28').Value = '  -6.95%  '
$ws.Range('D29').Value = '1.360'
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('D30').Value = '0.09621'
$ws.Range('E30').Value = '  -4.34%  '
$ws.Range('E31').Value = '  -5.73%  '
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('D33').Value = '4.067'
$ws.Range('D34').Value = '0.04633'
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').Value = '0.6950'
$ws.Range('E35').Value = '  -4.65%  '
$ws.Range('D36').Value = '1.078'
$ws.Range('E36').Value = '  -3.08%  '
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '2.702'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  -4.31%  '
$ws.Range('D40').Value = '6.325'
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = '2.501'
$ws.Range('E41').Value = '  -4.16%  '
$ws.Range('D42').Value = '71.02'
$ws.Range('E42').Value = '  -5.22%  '
$ws.Range('D43').Value = '0.8548'
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('D44').Value = '1.923'
$ws.Range('E44').Value = '  -2.42%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').Value = '102.70'
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '0.4135'
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('D48').Value = '973.35'
$ws.Range('E48').Value = '  +5.88%  '
$ws.Range('E49').Value = '  -3.92%  '
$ws.Range('D50').Value = '9.034'
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('E51').Value = '  -3.59%  '

# Restore original (unstyled) formatting so styles.xml / cell style indices
# remain unaffected, matching the source diff which only touches cell text.
foreach ($ref in $targets) {
    $ws.Range($ref).ClearFormats()
}
